$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-empty "resultado" / "profit" cells for already-existing rows
$ws.Range("G78").Value = "Acierto"
$ws.Range("H78").Value = 1.5

$ws.Range("G79").Value = "Acierto"
$ws.Range("H79").Value = 1.1

$ws.Range("G84").Value = "Fallo"
$ws.Range("H84").Value = -1

$ws.Range("G88").Value = "Acierto"
$ws.Range("H88").Value = 1

$ws.Range("G90").Value = "Acierto"
$ws.Range("H90").Value = 1.62

# Append two new result rows at the bottom of the tracker.
# "fecha" (column B) holds text like "2025-08-06" in this sheet (not a real
# date), so force the cell to Text format before writing the value -- and
# then clear the formatting again -- to stop Excel from auto-converting the
# literal into a date serial number.
$ws.Range("A91").Value = 14266328
$ws.Range("B91").NumberFormat = "@"
$ws.Range("B91").Value = "2025-08-06"
$ws.Range("B91").ClearFormats()
$ws.Range("C91").Value = "Victoria Mboko"
$ws.Range("D91").Value = "Elena Rybakina"
$ws.Range("E91").Value = "Gana Victoria Mboko"
$ws.Range("F91").Value = 4
$ws.Range("G91").Value = "Acierto"
$ws.Range("H91").Value = 3

$ws.Range("A92").Value = 14367423
$ws.Range("B92").NumberFormat = "@"
$ws.Range("B92").Value = "2025-08-06"
$ws.Range("B92").ClearFormats()
$ws.Range("C92").Value = "Aleksandar Vukic"
$ws.Range("D92").Value = "Patrick Kypson"
$ws.Range("E92").Value = "Gana Aleksandar Vukic"
$ws.Range("F92").Value = 2.62
$ws.Range("G92").Value = "Fallo"
$ws.Range("H92").Value = -1
